$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.917.32"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +1.74%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.248.08"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.04%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.06%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'318.70"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -0.09%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'100.37"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +1.60%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.573"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  -1.29%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.02%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.547"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -3.25%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'36.86"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.05%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.0830"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.56%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'7.54"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  -1.56%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.91%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.591.24"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.08%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'14.42"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +0.58%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.853"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  -0.90%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.251.89"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.86%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'43.769.66"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +1.56%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'13.50"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  -1.92%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'0.0₃0978"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.55%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  -1.73%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'65.37"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.31%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  -4.54%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'233.84"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -1.02%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').Value = "'  -6.14%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.08%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'10.72"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +6.57%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'38.75"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +5.10%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.18"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -1.53%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D31').Value = "'160.36"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.85%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'20.07"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -0.73%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').Value = "'  -2.54%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'2.68"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.02%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('B35').Value = "'Kaspa"
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = "'0.114"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +8.82%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('B36').Value = "'LidoDAOToken"
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = "'3.09"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  -6.53%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +5.47%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.118"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.71%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'16.48"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +16.16%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'3.67"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  -0.32%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'4.15"
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').Value = "'0.0314"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -1.47%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('E43').Value = "'  +0.01%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'1.769.02"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.08%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('B45').Value = "'ordi"
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = "'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = "'74.60"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +1.12%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = "'Algorand"
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = "'0.195"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -3.61%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'5.16"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  -1.99%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'80.91"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -3.26%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'103.73"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.99%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('B50').Value = "'Stacks"
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = "'1.66"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +3.14%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('B51').Value = "'MultiversX"
$ws.Range('B51').Style = 'Normal'
$ws.Range('C51').Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range('C51').Style = 'Normal'
$ws.Range('D51').Value = "'57.25"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -0.75%  "
$ws.Range('E51').Style = 'Normal'
